# edit.ps1 - Apply "feat: add 2022-Q3 data" change
# 1. Insert a new worksheet "2022-Q3" right after "总计", containing the latest
#    fund-holdings snapshot.
# 2. Insert a new summary row into "总计" for 2022-Q3 (existing rows shift down).

function Set-HeaderStyle($cell) {
    # Reproduce the bold / thin-border / center-top style used throughout
    # the workbook for header cells and the index column.
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

function Set-TextValue($cell, $val) {
    # Force the value to be stored as text (matches the source data, which
    # keeps numeric-looking strings like fund codes / percentages as text).
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the new "2022-Q3" worksheet right after "总计" (position 2)
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$ws = $wb.Worksheets.Add($null, $totalSheet)
$ws.Name = "2022-Q3"

    # Header row
    $c = $ws.Cells.Item(1, 2)
    Set-HeaderStyle $c
    $c.Value = "基金代码"
    $c = $ws.Cells.Item(1, 3)
    Set-HeaderStyle $c
    $c.Value = "基金名称"
    $c = $ws.Cells.Item(1, 4)
    Set-HeaderStyle $c
    $c.Value = "基金规模"
    $c = $ws.Cells.Item(1, 5)
    Set-HeaderStyle $c
    $c.Value = "股票总仓位"
    $c = $ws.Cells.Item(1, 6)
    Set-HeaderStyle $c
    $c.Value = "仓位占比"
    $c = $ws.Cells.Item(1, 7)
    Set-HeaderStyle $c
    $c.Value = "持有市值(亿元)"
    $c = $ws.Cells.Item(1, 8)
    Set-HeaderStyle $c
    $c.Value = "仓位排名"

    # Data rows
    # row 2
    $c = $ws.Cells.Item(2, 1)
    Set-HeaderStyle $c
    $c.Value = 0
    Set-TextValue $ws.Cells.Item(2, 2) "512880"
    Set-TextValue $ws.Cells.Item(2, 3) "国泰中证全指证券公司ETF"
    Set-TextValue $ws.Cells.Item(2, 4) "289.27"
    Set-TextValue $ws.Cells.Item(2, 5) "99.95"
    Set-TextValue $ws.Cells.Item(2, 6) "2.66"
    Set-TextValue $ws.Cells.Item(2, 7) "7.6946"
    $ws.Cells.Item(2, 8).Value = 10
    # row 3
    $c = $ws.Cells.Item(3, 1)
    Set-HeaderStyle $c
    $c.Value = 1
    Set-TextValue $ws.Cells.Item(3, 2) "512000"
    Set-TextValue $ws.Cells.Item(3, 3) "华宝中证全指证券公司ETF"
    Set-TextValue $ws.Cells.Item(3, 4) "215.91"
    Set-TextValue $ws.Cells.Item(3, 5) "99.89"
    Set-TextValue $ws.Cells.Item(3, 6) "2.66"
    Set-TextValue $ws.Cells.Item(3, 7) "5.7432"
    $ws.Cells.Item(3, 8).Value = 10
    # row 4
    $c = $ws.Cells.Item(4, 1)
    Set-HeaderStyle $c
    $c.Value = 2
    Set-TextValue $ws.Cells.Item(4, 2) "512900"
    Set-TextValue $ws.Cells.Item(4, 3) "南方中证全指证券公司ETF"
    Set-TextValue $ws.Cells.Item(4, 4) "78.74"
    Set-TextValue $ws.Cells.Item(4, 5) "99.98"
    Set-TextValue $ws.Cells.Item(4, 6) "2.66"
    Set-TextValue $ws.Cells.Item(4, 7) "2.0945"
    $ws.Cells.Item(4, 8).Value = 10
    # row 5
    $c = $ws.Cells.Item(5, 1)
    Set-HeaderStyle $c
    $c.Value = 3
    Set-TextValue $ws.Cells.Item(5, 2) "159841"
    Set-TextValue $ws.Cells.Item(5, 3) "天弘中证全指证券公司ETF"
    Set-TextValue $ws.Cells.Item(5, 4) "45.70"
    Set-TextValue $ws.Cells.Item(5, 5) "99.94"
    Set-TextValue $ws.Cells.Item(5, 6) "2.66"
    Set-TextValue $ws.Cells.Item(5, 7) "1.2156"
    $ws.Cells.Item(5, 8).Value = 10
    # row 6
    $c = $ws.Cells.Item(6, 1)
    Set-HeaderStyle $c
    $c.Value = 4
    Set-TextValue $ws.Cells.Item(6, 2) "161720"
    Set-TextValue $ws.Cells.Item(6, 3) "招商中证全指证券公司指数（LOF）A"
    Set-TextValue $ws.Cells.Item(6, 4) "22.28"
    Set-TextValue $ws.Cells.Item(6, 5) "94.50"
    Set-TextValue $ws.Cells.Item(6, 6) "2.51"
    Set-TextValue $ws.Cells.Item(6, 7) "0.5592"
    $ws.Cells.Item(6, 8).Value = 10
    # row 7
    $c = $ws.Cells.Item(7, 1)
    Set-HeaderStyle $c
    $c.Value = 5
    Set-TextValue $ws.Cells.Item(7, 2) "501016"
    Set-TextValue $ws.Cells.Item(7, 3) "国泰中证申万证券行业指数（LOF）A"
    Set-TextValue $ws.Cells.Item(7, 4) "18.10"
    Set-TextValue $ws.Cells.Item(7, 5) "93.42"
    Set-TextValue $ws.Cells.Item(7, 6) "2.52"
    Set-TextValue $ws.Cells.Item(7, 7) "0.4561"
    $ws.Cells.Item(7, 8).Value = 10
    # row 8
    $c = $ws.Cells.Item(8, 1)
    Set-HeaderStyle $c
    $c.Value = 6
    Set-TextValue $ws.Cells.Item(8, 2) "163113"
    Set-TextValue $ws.Cells.Item(8, 3) "申万菱信中证申万证券行业指数（LOF）A"
    Set-TextValue $ws.Cells.Item(8, 4) "16.20"
    Set-TextValue $ws.Cells.Item(8, 5) "93.19"
    Set-TextValue $ws.Cells.Item(8, 6) "2.51"
    Set-TextValue $ws.Cells.Item(8, 7) "0.4066"
    $ws.Cells.Item(8, 8).Value = 10
    # row 9
    $c = $ws.Cells.Item(9, 1)
    Set-HeaderStyle $c
    $c.Value = 7
    Set-TextValue $ws.Cells.Item(9, 2) "161027"
    Set-TextValue $ws.Cells.Item(9, 3) "富国中证全指证券公司指数A"
    Set-TextValue $ws.Cells.Item(9, 4) "12.47"
    Set-TextValue $ws.Cells.Item(9, 5) "94.32"
    Set-TextValue $ws.Cells.Item(9, 6) "2.51"
    Set-TextValue $ws.Cells.Item(9, 7) "0.3130"
    $ws.Cells.Item(9, 8).Value = 10
    # row 10
    $c = $ws.Cells.Item(10, 1)
    Set-HeaderStyle $c
    $c.Value = 8
    Set-TextValue $ws.Cells.Item(10, 2) "502010"
    Set-TextValue $ws.Cells.Item(10, 3) "易方达证券公司指数（LOF）A"
    Set-TextValue $ws.Cells.Item(10, 4) "12.08"
    Set-TextValue $ws.Cells.Item(10, 5) "94.58"
    Set-TextValue $ws.Cells.Item(10, 6) "2.51"
    Set-TextValue $ws.Cells.Item(10, 7) "0.3032"
    $ws.Cells.Item(10, 8).Value = 10
    # row 11
    $c = $ws.Cells.Item(11, 1)
    Set-HeaderStyle $c
    $c.Value = 9
    Set-TextValue $ws.Cells.Item(11, 2) "160633"
    Set-TextValue $ws.Cells.Item(11, 3) "鹏华中证全指证券公司指数（LOF）A"
    Set-TextValue $ws.Cells.Item(11, 4) "11.96"
    Set-TextValue $ws.Cells.Item(11, 5) "94.06"
    Set-TextValue $ws.Cells.Item(11, 6) "2.50"
    Set-TextValue $ws.Cells.Item(11, 7) "0.2990"
    $ws.Cells.Item(11, 8).Value = 10
    # row 12
    $c = $ws.Cells.Item(12, 1)
    Set-HeaderStyle $c
    $c.Value = 10
    Set-TextValue $ws.Cells.Item(12, 2) "501048"
    Set-TextValue $ws.Cells.Item(12, 3) "汇添富中证全指证券公司指数（LOF）C"
    Set-TextValue $ws.Cells.Item(12, 4) "9.31"
    Set-TextValue $ws.Cells.Item(12, 5) "93.53"
    Set-TextValue $ws.Cells.Item(12, 6) "2.50"
    Set-TextValue $ws.Cells.Item(12, 7) "0.2328"
    $ws.Cells.Item(12, 8).Value = 10
    # row 13
    $c = $ws.Cells.Item(13, 1)
    Set-HeaderStyle $c
    $c.Value = 11
    Set-TextValue $ws.Cells.Item(13, 2) "515010"
    Set-TextValue $ws.Cells.Item(13, 3) "华夏中证全指证券公司ETF"
    Set-TextValue $ws.Cells.Item(13, 4) "8.47"
    Set-TextValue $ws.Cells.Item(13, 5) "99.66"
    Set-TextValue $ws.Cells.Item(13, 6) "2.64"
    Set-TextValue $ws.Cells.Item(13, 7) "0.2236"
    $ws.Cells.Item(13, 8).Value = 10
    # row 14
    $c = $ws.Cells.Item(14, 1)
    Set-HeaderStyle $c
    $c.Value = 12
    Set-TextValue $ws.Cells.Item(14, 2) "160516"
    Set-TextValue $ws.Cells.Item(14, 3) "博时中证全指证券公司指数"
    Set-TextValue $ws.Cells.Item(14, 4) "7.00"
    Set-TextValue $ws.Cells.Item(14, 5) "93.76"
    Set-TextValue $ws.Cells.Item(14, 6) "2.49"
    Set-TextValue $ws.Cells.Item(14, 7) "0.1743"
    $ws.Cells.Item(14, 8).Value = 10
    # row 15
    $c = $ws.Cells.Item(15, 1)
    Set-HeaderStyle $c
    $c.Value = 13
    Set-TextValue $ws.Cells.Item(15, 2) "012044"
    Set-TextValue $ws.Cells.Item(15, 3) "鹏华中证全指证券公司指数（LOF）C"
    Set-TextValue $ws.Cells.Item(15, 4) "5.89"
    Set-TextValue $ws.Cells.Item(15, 5) "94.06"
    Set-TextValue $ws.Cells.Item(15, 6) "2.50"
    Set-TextValue $ws.Cells.Item(15, 7) "0.1472"
    $ws.Cells.Item(15, 8).Value = 10
    # row 16
    $c = $ws.Cells.Item(16, 1)
    Set-HeaderStyle $c
    $c.Value = 14
    Set-TextValue $ws.Cells.Item(16, 2) "501047"
    Set-TextValue $ws.Cells.Item(16, 3) "汇添富中证全指证券公司指数（LOF）A"
    Set-TextValue $ws.Cells.Item(16, 4) "5.58"
    Set-TextValue $ws.Cells.Item(16, 5) "93.53"
    Set-TextValue $ws.Cells.Item(16, 6) "2.50"
    Set-TextValue $ws.Cells.Item(16, 7) "0.1395"
    $ws.Cells.Item(16, 8).Value = 10
    # row 17
    $c = $ws.Cells.Item(17, 1)
    Set-HeaderStyle $c
    $c.Value = 15
    Set-TextValue $ws.Cells.Item(17, 2) "159842"
    Set-TextValue $ws.Cells.Item(17, 3) "银华中证全指证券公司ETF"
    Set-TextValue $ws.Cells.Item(17, 4) "5.01"
    Set-TextValue $ws.Cells.Item(17, 5) "98.00"
    Set-TextValue $ws.Cells.Item(17, 6) "2.60"
    Set-TextValue $ws.Cells.Item(17, 7) "0.1303"
    $ws.Cells.Item(17, 8).Value = 10
    # row 18
    $c = $ws.Cells.Item(18, 1)
    Set-HeaderStyle $c
    $c.Value = 16
    Set-TextValue $ws.Cells.Item(18, 2) "515560"
    Set-TextValue $ws.Cells.Item(18, 3) "建信中证全指证券公司ETF"
    Set-TextValue $ws.Cells.Item(18, 4) "3.98"
    Set-TextValue $ws.Cells.Item(18, 5) "98.63"
    Set-TextValue $ws.Cells.Item(18, 6) "2.63"
    Set-TextValue $ws.Cells.Item(18, 7) "0.1047"
    $ws.Cells.Item(18, 8).Value = 10
    # row 19
    $c = $ws.Cells.Item(19, 1)
    Set-HeaderStyle $c
    $c.Value = 17
    Set-TextValue $ws.Cells.Item(19, 2) "502053"
    Set-TextValue $ws.Cells.Item(19, 3) "长盛中证全指证券公司指数（LOF）"
    Set-TextValue $ws.Cells.Item(19, 4) "3.93"
    Set-TextValue $ws.Cells.Item(19, 5) "93.15"
    Set-TextValue $ws.Cells.Item(19, 6) "2.49"
    Set-TextValue $ws.Cells.Item(19, 7) "0.0979"
    $ws.Cells.Item(19, 8).Value = 10
    # row 20
    $c = $ws.Cells.Item(20, 1)
    Set-HeaderStyle $c
    $c.Value = 18
    Set-TextValue $ws.Cells.Item(20, 2) "160419"
    Set-TextValue $ws.Cells.Item(20, 3) "华安中证证券公司A"
    Set-TextValue $ws.Cells.Item(20, 4) "3.88"
    Set-TextValue $ws.Cells.Item(20, 5) "94.43"
    Set-TextValue $ws.Cells.Item(20, 6) "2.50"
    Set-TextValue $ws.Cells.Item(20, 7) "0.0970"
    $ws.Cells.Item(20, 8).Value = 10
    # row 21
    $c = $ws.Cells.Item(21, 1)
    Set-HeaderStyle $c
    $c.Value = 19
    Set-TextValue $ws.Cells.Item(21, 2) "012874"
    Set-TextValue $ws.Cells.Item(21, 3) "易方达证券公司指数（LOF）C"
    Set-TextValue $ws.Cells.Item(21, 4) "2.22"
    Set-TextValue $ws.Cells.Item(21, 5) "94.58"
    Set-TextValue $ws.Cells.Item(21, 6) "2.51"
    Set-TextValue $ws.Cells.Item(21, 7) "0.0557"
    $ws.Cells.Item(21, 8).Value = 10
    # row 22
    $c = $ws.Cells.Item(22, 1)
    Set-HeaderStyle $c
    $c.Value = 20
    Set-TextValue $ws.Cells.Item(22, 2) "512570"
    Set-TextValue $ws.Cells.Item(22, 3) "易方达中证全指证券公司ETF"
    Set-TextValue $ws.Cells.Item(22, 4) "1.92"
    Set-TextValue $ws.Cells.Item(22, 5) "98.99"
    Set-TextValue $ws.Cells.Item(22, 6) "2.63"
    Set-TextValue $ws.Cells.Item(22, 7) "0.0505"
    $ws.Cells.Item(22, 8).Value = 10
    # row 23
    $c = $ws.Cells.Item(23, 1)
    Set-HeaderStyle $c
    $c.Value = 21
    Set-TextValue $ws.Cells.Item(23, 2) "515850"
    Set-TextValue $ws.Cells.Item(23, 3) "富国中证全指证券公司ETF"
    Set-TextValue $ws.Cells.Item(23, 4) "1.67"
    Set-TextValue $ws.Cells.Item(23, 5) "99.74"
    Set-TextValue $ws.Cells.Item(23, 6) "2.67"
    Set-TextValue $ws.Cells.Item(23, 7) "0.0446"
    $ws.Cells.Item(23, 8).Value = 10
    # row 24
    $c = $ws.Cells.Item(24, 1)
    Set-HeaderStyle $c
    $c.Value = 22
    Set-TextValue $ws.Cells.Item(24, 2) "159848"
    Set-TextValue $ws.Cells.Item(24, 3) "国联安中证全指证券公司ETF"
    Set-TextValue $ws.Cells.Item(24, 4) "0.96"
    Set-TextValue $ws.Cells.Item(24, 5) "96.87"
    Set-TextValue $ws.Cells.Item(24, 6) "2.58"
    Set-TextValue $ws.Cells.Item(24, 7) "0.0248"
    $ws.Cells.Item(24, 8).Value = 10
    # row 25
    $c = $ws.Cells.Item(25, 1)
    Set-HeaderStyle $c
    $c.Value = 23
    Set-TextValue $ws.Cells.Item(25, 2) "516730"
    Set-TextValue $ws.Cells.Item(25, 3) "浦银安盛中证证券公司30ETF"
    Set-TextValue $ws.Cells.Item(25, 4) "0.60"
    Set-TextValue $ws.Cells.Item(25, 5) "97.43"
    Set-TextValue $ws.Cells.Item(25, 6) "3.27"
    Set-TextValue $ws.Cells.Item(25, 7) "0.0196"
    $ws.Cells.Item(25, 8).Value = 10
    # row 26
    $c = $ws.Cells.Item(26, 1)
    Set-HeaderStyle $c
    $c.Value = 24
    Set-TextValue $ws.Cells.Item(26, 2) "013712"
    Set-TextValue $ws.Cells.Item(26, 3) "方正富邦鑫益一年定期开放混合A"
    Set-TextValue $ws.Cells.Item(26, 4) "2.04"
    Set-TextValue $ws.Cells.Item(26, 5) "32.96"
    Set-TextValue $ws.Cells.Item(26, 6) "0.92"
    Set-TextValue $ws.Cells.Item(26, 7) "0.0188"
    $ws.Cells.Item(26, 8).Value = 9
    # row 27
    $c = $ws.Cells.Item(27, 1)
    Set-HeaderStyle $c
    $c.Value = 25
    Set-TextValue $ws.Cells.Item(27, 2) "013276"
    Set-TextValue $ws.Cells.Item(27, 3) "富国中证全指证券公司指数C"
    Set-TextValue $ws.Cells.Item(27, 4) "0.57"
    Set-TextValue $ws.Cells.Item(27, 5) "94.32"
    Set-TextValue $ws.Cells.Item(27, 6) "2.51"
    Set-TextValue $ws.Cells.Item(27, 7) "0.0143"
    $ws.Cells.Item(27, 8).Value = 10
    # row 28
    $c = $ws.Cells.Item(28, 1)
    Set-HeaderStyle $c
    $c.Value = 26
    Set-TextValue $ws.Cells.Item(28, 2) "013597"
    Set-TextValue $ws.Cells.Item(28, 3) "招商中证全指证券公司指数（LOF）C"
    Set-TextValue $ws.Cells.Item(28, 4) "0.39"
    Set-TextValue $ws.Cells.Item(28, 5) "94.50"
    Set-TextValue $ws.Cells.Item(28, 6) "2.51"
    Set-TextValue $ws.Cells.Item(28, 7) "0.0098"
    $ws.Cells.Item(28, 8).Value = 10
    # row 29
    $c = $ws.Cells.Item(29, 1)
    Set-HeaderStyle $c
    $c.Value = 27
    Set-TextValue $ws.Cells.Item(29, 2) "516200"
    Set-TextValue $ws.Cells.Item(29, 3) "华安中证全指证券公司ETF"
    Set-TextValue $ws.Cells.Item(29, 4) "0.32"
    Set-TextValue $ws.Cells.Item(29, 5) "97.22"
    Set-TextValue $ws.Cells.Item(29, 6) "2.60"
    Set-TextValue $ws.Cells.Item(29, 7) "0.0083"
    $ws.Cells.Item(29, 8).Value = 10
    # row 30
    $c = $ws.Cells.Item(30, 1)
    Set-HeaderStyle $c
    $c.Value = 28
    Set-TextValue $ws.Cells.Item(30, 2) "516980"
    Set-TextValue $ws.Cells.Item(30, 3) "华富中证证券公司先锋策略ETF"
    Set-TextValue $ws.Cells.Item(30, 4) "0.28"
    Set-TextValue $ws.Cells.Item(30, 5) "99.02"
    Set-TextValue $ws.Cells.Item(30, 6) "2.60"
    Set-TextValue $ws.Cells.Item(30, 7) "0.0073"
    $ws.Cells.Item(30, 8).Value = 9
    # row 31
    $c = $ws.Cells.Item(31, 1)
    Set-HeaderStyle $c
    $c.Value = 29
    Set-TextValue $ws.Cells.Item(31, 2) "014984"
    Set-TextValue $ws.Cells.Item(31, 3) "华安中证证券公司C"
    Set-TextValue $ws.Cells.Item(31, 4) "0.18"
    Set-TextValue $ws.Cells.Item(31, 5) "94.43"
    Set-TextValue $ws.Cells.Item(31, 6) "2.50"
    Set-TextValue $ws.Cells.Item(31, 7) "0.0045"
    $ws.Cells.Item(31, 8).Value = 10
    # row 32
    $c = $ws.Cells.Item(32, 1)
    Set-HeaderStyle $c
    $c.Value = 30
    Set-TextValue $ws.Cells.Item(32, 2) "015178"
    Set-TextValue $ws.Cells.Item(32, 3) "申万菱信中证申万证券行业指数（LOF）C"
    Set-TextValue $ws.Cells.Item(32, 4) "0.08"
    Set-TextValue $ws.Cells.Item(32, 5) "93.19"
    Set-TextValue $ws.Cells.Item(32, 6) "2.51"
    Set-TextValue $ws.Cells.Item(32, 7) "0.0020"
    $ws.Cells.Item(32, 8).Value = 10
    # row 33
    $c = $ws.Cells.Item(33, 1)
    Set-HeaderStyle $c
    $c.Value = 31
    Set-TextValue $ws.Cells.Item(33, 2) "013713"
    Set-TextValue $ws.Cells.Item(33, 3) "方正富邦鑫益一年定期开放混合C"
    Set-TextValue $ws.Cells.Item(33, 4) "0.06"
    Set-TextValue $ws.Cells.Item(33, 5) "32.96"
    Set-TextValue $ws.Cells.Item(33, 6) "0.92"
    Set-TextValue $ws.Cells.Item(33, 7) "0.0006"
    $ws.Cells.Item(33, 8).Value = 9
    # row 34
    $c = $ws.Cells.Item(34, 1)
    Set-HeaderStyle $c
    $c.Value = 32
    Set-TextValue $ws.Cells.Item(34, 2) "015598"
    Set-TextValue $ws.Cells.Item(34, 3) "国泰中证申万证券行业指数（LOF）C"
    Set-TextValue $ws.Cells.Item(34, 4) "0.01"
    Set-TextValue $ws.Cells.Item(34, 5) "93.42"
    Set-TextValue $ws.Cells.Item(34, 6) "2.52"
    Set-TextValue $ws.Cells.Item(34, 7) "0.0003"
    $ws.Cells.Item(34, 8).Value = 10

# ---------------------------------------------------------------------------
# Step 2: insert the 2022-Q3 summary row into "总计" (existing rows shift down)
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

$cA = $totalSheet.Cells.Item(2, 1)
Set-HeaderStyle $cA
$cA.Value = 0

Set-TextValue $totalSheet.Cells.Item(2, 2) "2022-Q3"

$totalSheet.Cells.Item(2, 3).Value = 33
$totalSheet.Cells.Item(2, 4).Value = 20.69

# Keep the originally-active sheet ("总计") selected, since the source diff
# does not indicate any change to the active tab.
$totalSheet.Activate()

Write-Host "Done applying 2022-Q3 update"
